$d = $word.ActiveDocument

$b = $d.Bookmarks.Item("_GoBack")
$b.Delete()

$p6 = $d.Paragraphs.Item(6)
$r6 = $p6.Range
$xml6 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Diec Du </w:t></w:r><w:r><w:t>TRAN</w:t></w:r><w:r><w:t xml:space="preserve"> changing things </w:t></w:r><w:r><w:t>up too.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r6.InsertXML($xml6)

$p7 = $d.Paragraphs.Item(7)
$r7 = $p7.Range
$xml7 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Then complete and close this document.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r7.InsertXML($xml7)

# Delete paragraph 8's own mark (the trailing empty original paragraph) by
# deleting a zero-width range expanded to include just its 1-char mark,
# then see which attrs survive.
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$lr = $last.Range
Write-Output ("last: " + $lr.Start + "-" + $lr.End)
$delRange = $d.Range($lr.Start - 1, $lr.End)
Write-Output ("delRange=[" + $delRange.Text + "] " + $delRange.Start + "-" + $delRange.End)
$delRange.Delete()

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output ("Para " + $i + ": Start=" + $p.Range.Start + " End=" + $p.Range.End + " Text=[" + $p.Range.Text + "]")
}
